# Applies the ECOLAND_20250708_cycle3 edit:
#  - Raw_Annotations!E2:E31 all changed to 300 (traffic-count normalisation window),
#    which cascades through the shared J/K formulas automatically.
#  - Active sheet/tab moves from Aggregates -> Raw_Annotations, with the
#    selection on Raw_Annotations now covering E2:E31 (the column that changed).

$wb = $excel.ActiveWorkbook

$wsRaw = $wb.Worksheets.Item("Raw_Annotations")

# --- Update the raw annotation interval column (E2:E31) to 300 ---
$wsRaw.Range("E2:E31").Value = 300

# --- Move the active tab / selection to Raw_Annotations, matching the new
#     selection rectangle E2:E31 left over from the bulk edit above ---
$wsRaw.Activate()
[void]$wsRaw.Range("E2:E31").Select()

Write-Host "Updated E2:E31 to 300 on Raw_Annotations and re-activated the sheet"
